$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, so Excel stores the literal string (matching the source data)
# instead of silently converting to a numeric value.
$textCells = @("D5", "D6", "D8", "D11", "D18", "D19", "D25", "D26", "D27", "D36", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "34.561.59"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.810.36"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "225.92"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "0.597"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "37.53"
$ws.Range("E8").Value = "  +7.48%  "
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "0.0969"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "2.071.09"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "1.813.61"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "34.526.86"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "68.65"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "244.17"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  +4.36%  "
$ws.Range("D25").Value = "172.13"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "7.86"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("D27").Value = "17.38"
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").Value = "1.365.70"
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("D36").Value = "0.655"
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").Value = "0.0188"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("E40").Value = "  +8.25%  "
$ws.Range("D41").Value = "2.42"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").Value = "80.97"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "0.940"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").Value = "2.78"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").Value = "13.85"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("D46").Value = "0.0500"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").Value = "1.971.53"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "103.08"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("E51").Value = "  -6.88%  "
